# Adding new stats for barriers and data.csv
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Block 1 (rows 1-5): Mammogram x Breast-exam odds-ratio table
# ---------------------------------------------------------------
$ws.Range("Q1").Value = "No breast exam"
$ws.Range("R1").Value = "Breast exam"
$ws.Range("S1").Value = "NB/B"

$ws.Range("P2").Value = "Hispanic"
$ws.Range("Q2").Value = 10
$ws.Range("R2").Value = 21
$ws.Range("S2").Formula = "=Q2/R2"

$ws.Range("P3").Value = "White"
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = 12
$ws.Range("S3").Formula = "=Q3/R3"

$ws.Range("R5").Value = "Odds Ratio="
$ws.Range("S5").Formula = "=S2/S3"

# ---------------------------------------------------------------
# Block 2 (rows 7-11): Non-compliant/Compliant x Breast-exam table
# ---------------------------------------------------------------
$ws.Range("Q7").Value = "Non-compliant"
$ws.Range("R7").Value = "Compliant"

$ws.Range("P8").Value = "Hispanic"
$ws.Range("Q8").Value = 7
$ws.Range("R8").Value = 10
$ws.Range("S8").Formula = "=Q8/R8"

$ws.Range("P9").Value = "White"
$ws.Range("Q9").Value = 1
$ws.Range("R9").Value = 11
$ws.Range("S9").Formula = "=Q9/R9"

$ws.Range("R11").Value = "Odds Ratio="
$ws.Range("S11").Formula = "=S8/S9"

# ---------------------------------------------------------------
# Block 3 (rows 13-18): INCOME x Breast-exam table
# ---------------------------------------------------------------
$ws.Range("P13").Value = "INCOME"

$ws.Range("Q14").Value = "No BE"
$ws.Range("R14").Value = "BE"

$ws.Range("P15").Value = "<40k"
$ws.Range("Q15").Value = 10
$ws.Range("R15").Value = 27
$ws.Range("S15").Formula = "=Q15/R15"

$ws.Range("P16").Value = ">40k"
$ws.Range("Q16").Value = 3
$ws.Range("R16").Value = 19
$ws.Range("S16").Formula = "=Q16/R16"

$ws.Range("R18").Value = "Odds Ratio="
$ws.Range("S18").Formula = "=S15/S16"

# ---------------------------------------------------------------
# Block 4 (rows 20-24): INCOME x Non-compliant/Compliant table
# ---------------------------------------------------------------
$ws.Range("Q20").Value = "NC"
$ws.Range("R20").Value = "Compliant"

$ws.Range("P21").Value = "<40k"
$ws.Range("Q21").Value = 7
$ws.Range("R21").Value = 16
$ws.Range("S21").Formula = "=Q21/R21"

$ws.Range("P22").Value = ">40k"
$ws.Range("Q22").Value = 2
$ws.Range("R22").Value = 14
$ws.Range("S22").Formula = "=Q22/R22"

$ws.Range("R24").Value = "Odds Ratio="
$ws.Range("S24").Formula = "=S21/S22"

# ---------------------------------------------------------------
# Block 5 (rows 25-32): EDUCATION barriers frequency/percentage
#                        table + EDUCATION x Breast-exam table
# ---------------------------------------------------------------
$ws.Range("K25").Value = "Barriers"
$ws.Range("K25").Font.Bold = $true
$ws.Range("M25").Value = "Percentage"
$ws.Range("P25").Value = "EDUCATION"

$ws.Range("K26").Value = 1
$ws.Range("L26").Value = 11
$ws.Range("M26").Formula = "=L26/48%"
$ws.Range("Q26").Value = "No BE"
$ws.Range("R26").Value = "BE"

$ws.Range("K27").Value = 2
$ws.Range("L27").Value = 1
$ws.Range("P27").Value = "<=HS"
$ws.Range("Q27").Value = 5
$ws.Range("R27").Value = 16
$ws.Range("S27").Formula = "=Q27/R27"

$ws.Range("K28").Value = 3
$ws.Range("L28").Value = 5
$ws.Range("P28").Value = ">HS"
$ws.Range("Q28").Value = 8
$ws.Range("R28").Value = 30
$ws.Range("S28").Formula = "=Q28/R28"

$ws.Range("K29").Value = 4
$ws.Range("L29").Value = 1

$ws.Range("K30").Value = 5
$ws.Range("L30").Value = 19
$ws.Range("R30").Value = "Odds Ratio="
$ws.Range("S30").Formula = "=S27/S28"

$ws.Range("K31").Value = 6
$ws.Range("L31").Value = 11

# M27:M31 share one relative formula anchored at M27 (matches how Excel
# fills a formula down a multi-cell selection).
$ws.Range("M27:M31").Formula = "=L27/48%"

$ws.Range("K32").Value = "Total"
$ws.Range("L32").Formula = "=SUM(L26:L31)"
$ws.Range("Q32").Value = "NC"
$ws.Range("R32").Value = "Compliant"

# ---------------------------------------------------------------
# Block 6 (rows 33-36): EDUCATION x Non-compliant/Compliant table
# ---------------------------------------------------------------
$ws.Range("P33").Value = "<=HS"
$ws.Range("Q33").Value = 5
$ws.Range("R33").Value = 7
$ws.Range("S33").Formula = "=Q33/R33"

$ws.Range("P34").Value = ">HS"
$ws.Range("Q34").Value = 4
$ws.Range("R34").Value = 23
$ws.Range("S34").Formula = "=Q34/R34"

$ws.Range("R36").Value = "Odds Ratio="
$ws.Range("S36").Formula = "=S33/S34"

# ---------------------------------------------------------------
# Block 7 (rows 38-45): Barriers frequency/percentage tables
# (A/B/C = data.csv sample, F/G/H = second sample)
# ---------------------------------------------------------------
$ws.Range("A38").Value = "Barriers"
$ws.Range("A38").Font.Bold = $true
$ws.Range("C38").Value = "Percentage"
$ws.Range("F38").Value = "Barriers"
$ws.Range("F38").Font.Bold = $true
$ws.Range("H38").Value = "Percentage"

$ws.Range("A39").Value = 1
$ws.Range("B39").Value = 12
$ws.Range("C39").Formula = "=B39/B45 %"
$ws.Range("F39").Value = 1
$ws.Range("G39").Value = 21
$ws.Range("H39").Formula = "=G39/52 %"

$ws.Range("A40").Value = 2
$ws.Range("B40").Value = 2
$ws.Range("C40").Formula = "=B40/B45%"
$ws.Range("F40").Value = 2
$ws.Range("G40").Value = 3

$ws.Range("A41").Value = 3
$ws.Range("B41").Value = 3
$ws.Range("C41").Formula = "=B41/B45%"
$ws.Range("F41").Value = 3
$ws.Range("G41").Value = 12

$ws.Range("A42").Value = 4
$ws.Range("B42").Value = 1
$ws.Range("C42").Formula = "=B42/B45%"
$ws.Range("F42").Value = 4

$ws.Range("A43").Value = 5
$ws.Range("B43").Value = 5
$ws.Range("C43").Formula = "=B43/B45%"
$ws.Range("F43").Value = 5
$ws.Range("G43").Value = 5

$ws.Range("A44").Value = 6
$ws.Range("B44").Value = 6
$ws.Range("C44").Formula = "=B44/B45%"
$ws.Range("F44").Value = 6
$ws.Range("G44").Value = 11

# H40:H44 share one relative formula anchored at H40.
$ws.Range("H40:H44").Formula = "=G40/52 %"

$ws.Range("A45").Value = "Total"
$ws.Range("B45").Formula = "=SUM(B39:B44)"
$ws.Range("F45").Value = "Total"
$ws.Range("G45").Formula = "=SUM(G39:G44)"

# ---------------------------------------------------------------
# Column width for the new Q column + final selection
# ---------------------------------------------------------------
$ws.Columns("Q").ColumnWidth = 14.75

$ws.Range("M32").Select()
